# Apply updated crypto price/volume figures to worksheet (rows 2-51)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''66.979.77'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -0.95%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '''2.601.89'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -0.71%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = '''  -0.07%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''590.20'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -2.05%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''149.94'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -2.79%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = '''  -0.02%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = '''  -0.59%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = '''2.600.46'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -0.69%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''0.130'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '''  +2.32%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = '''  -0.12%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = '''5.15'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  -1.74%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = '''  -3.16%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '''27.18'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -2.81%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''3.070.59'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = '''0.0000182'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -3.08%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''66.854.47'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -1.16%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '''2.600.41'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  -0.80%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = '''  +0.10%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = '''11.03'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  -2.23%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''7.36'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -3.80%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = '''4.30'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -0.30%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = '''  -4.77%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = '''  -3.27%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''72.96'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  +9.75%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = '''  +0.03%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = '''9.89'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  -1.02%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = '''  -0.45%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''582.39'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  +0.00%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = '''  -1.13%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''0.0₃0989'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  -6.33%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = '''  -5.47%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("E34").Value = '''  -3.14%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = '''0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '''  -0.06%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = '''  -4.84%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = '''  -3.34%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''156.09'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '''  -1.02%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '''18.96'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '''  -2.54%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = '''  -1.70%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''1.85'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  +0.03%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''5.21'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  -3.59%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = '''  -5.04%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = '''  +3.96%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = '''  -0.06%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = '''152.83'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  -2.71%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''0.0₆0286'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  -2.34%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''3.66'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -2.86%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = '''  -3.62%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = '''  -1.64%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = '''21.43'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  +1.83%  '
$ws.Range("E51").Style = "Normal"
